$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '''42.612.27'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Formula = '''  -0.71%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Formula = '''2.265.06'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Formula = '''  -0.45%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').Formula = '''  -0.15%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Formula = '''250.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Formula = '''  +0.03%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('E6').Formula = '''  +0.71%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').Formula = '''75.67'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Formula = '''  +5.48%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('E8').Formula = '''  -0.02%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').Formula = '''0.640'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Formula = '''  -3.91%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').Formula = '''39.86'
$ws.Range('D10').Style = 'Normal'

$ws.Range('D11').Formula = '''0.0965'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Formula = '''  -0.84%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').Formula = '''7.29'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Formula = '''  -2.78%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').Formula = '''0.105'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Formula = '''  +1.07%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').Formula = '''2.603.98'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Formula = '''  -0.27%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').Formula = '''15.02'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Formula = '''  +0.68%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').Formula = '''0.862'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Formula = '''  -2.83%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Formula = '''2.270.35'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Formula = '''  -1.55%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Formula = '''42.513.72'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Formula = '''  -0.86%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('E19').Formula = '''  -1.20%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('E20').Formula = '''  -2.21%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').Formula = '''72.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Formula = '''  -1.67%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').Formula = '''233.91'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Formula = '''  -0.97%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('E23').Formula = '''  +1.39%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('E24').Formula = '''  +0.07%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').Formula = '''3.76'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Formula = '''  -5.14%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').Formula = '''11.24'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Formula = '''  -1.84%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('E27').Formula = '''  -2.13%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').Formula = '''2.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Formula = '''  -0.40%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').Formula = '''167.35'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Formula = '''  -0.24%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').Formula = '''20.90'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Formula = '''  -0.69%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').Formula = '''6.50'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Formula = '''  -1.67%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('D32').Formula = '''0.0858'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Formula = '''  +5.81%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('E33').Formula = '''  -2.98%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').Formula = '''31.45'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Formula = '''  +0.88%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('E35').Formula = '''  +0.73%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').Formula = '''4.54'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Formula = '''  +1.34%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('E38').Formula = '''  -3.78%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').Formula = '''13.71'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Formula = '''  +9.09%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('E40').Formula = '''  -3.41%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('D41').Formula = '''5.85'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Formula = '''  +0.17%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('E42').Formula = '''  +0.55%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').Formula = '''61.43'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Formula = '''  -1.28%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Formula = '''106.65'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Formula = '''  +11.79%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Formula = '''8.84'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Formula = '''  -4.61%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').Formula = '''4.71'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Formula = '''  -2.96%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('E47').Formula = '''  -1.91%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').Formula = '''0.998'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Formula = '''  -0.30%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('E49').Formula = '''  -2.77%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('E50').Formula = '''  -2.63%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('E51').Formula = '''  -2.17%  '
$ws.Range('E51').Style = 'Normal'
